$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 2 (current data row), shifting existing rows down.
$ws.Rows.Item(2).Insert()

# Excel's row insert copies formatting down from the row above (the bold header row).
# Reset the new row 2 to the plain (unstyled) formatting used by the other data rows,
# then restore the date number format on column D to match the other date cells.
$ws.Range("A2:R2").ClearFormats()
$ws.Range("D2").NumberFormat = $ws.Range("D3").NumberFormat

# New row 2 gets the same layout/values as the (now shifted) row 3, except for the updated fields.
$ws.Range("A2").Value = 10
$ws.Range("B2").Value = "Vega Modelo de Temuco"
$ws.Range("C2").Value = "La Araucanía"
$ws.Range("D2").Value = 44473
$ws.Range("E2").Value = 9
$ws.Range("F2").Value = 100112042
$ws.Range("G2").Value = "Locoto"
$ws.Range("H2").Value = "Sin especificar"
$ws.Range("I2").Value = "Primera"
$ws.Range("J2").Value = 140
$ws.Range("K2").Value = 1600
$ws.Range("L2").Value = 1600
$ws.Range("M2").Value = 1600
$ws.Range("N2").Value = "$/kilo"
$ws.Range("O2").Value = "Región de Arica y Parinacota"
$ws.Range("P2").Value = 1600
$ws.Range("Q2").Value = 1
$ws.Range("R2").Value = "Hortaliza"
